{"js": "// Highlight part of three requirement bullets in the \"Exam Shopping List\"\n// description document in yellow, by splitting the single run that holds\n// each sentence into multiple runs (matching Word's own behaviour when a\n// user selects a sub-string and applies a highlight color).\n//\n// Targets (exact paragraph text before edit):\n//   1. \"Name length must be between 3 and 20 characters (inclusive 3 and 20).\"\n//        -> highlight \"3 and 20 characters (inclusive 3 and 20).\"\n//   2. \"Description min length must be minimum 5(inclusive) characters\"\n//        -> highlight \"5(inclusive)\" only (leaving \" characters\" un-highlighted)\n//   3. \"Category cannot be null.\"\n//        -> highlight \"cannot be null.\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\n// Map exact paragraph text -> the substring(s) that must be highlighted.\n// Using the *whole* paragraph text as the lookup key guarantees we touch\n// only the intended bullet even though similar sentences (\"Username length\n// must be...\"/\"Password length must be...\") exist elsewhere in the doc.\nconst targets = [\n  {\n    paragraphText:\n      \"Name length must be between 3 and 20 characters (inclusive 3 and 20).\",\n    highlight: \"3 and 20 characters (inclusive 3 and 20).\",\n  },\n  {\n    paragraphText:\n      \"Description min length must be minimum 5(inclusive) characters\",\n    highlight: \"5(inclusive)\",\n  },\n  {\n    paragraphText: \"Category cannot be null.\",\n    highlight: \"cannot be null.\",\n  },\n];\n\nfor (const target of targets) {\n  const match = paragraphs.items.find((p) => p.text === target.paragraphText);\n  if (!match) {\n    continue;\n  }\n  const paragraphRange = match.getRange();\n  const found = paragraphRange.search(target.highlight, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length > 0) {\n    // There is exactly one occurrence within this paragraph; highlighting it\n    // causes Word to split the original run into separate runs so only the\n    // matched text carries the new formatting.\n    found.items[0].font.highlightColor = \"Yellow\";\n  }\n}\n\nawait context.sync();\n", "ps1": "# Highlight part of three requirement bullets in the \"Exam Shopping List\"\n# description document in yellow, by splitting the single run that holds\n# each sentence into multiple runs (matching Word's own behaviour when a\n# user selects a sub-string and applies a highlight color).\n#\n# Targets (exact paragraph text before edit):\n#   1. \"Name length must be between 3 and 20 characters (inclusive 3 and 20).\"\n#        -> highlight \"3 and 20 characters (inclusive 3 and 20).\"\n#   2. \"Description min length must be minimum 5(inclusive) characters\"\n#        -> highlight \"5(inclusive)\" only (leaving \" characters\" un-highlighted)\n#   3. \"Category cannot be null.\"\n#        -> highlight \"cannot be null.\"\n\n$d = $word.ActiveDocument\n\n# Map exact paragraph text -> the substring that must be highlighted.\n# Using the *whole* paragraph text as the lookup key guarantees we touch\n# only the intended bullet even though similar sentences (\"Username length\n# must be...\"/\"Password length must be...\") exist elsewhere in the doc.\n$targets = @(\n    @{ Paragraph = \"Name length must be between 3 and 20 characters (inclusive 3 and 20).\"; Highlight = \"3 and 20 characters (inclusive 3 and 20).\" },\n    @{ Paragraph = \"Description min length must be minimum 5(inclusive) characters\"; Highlight = \"5(inclusive)\" },\n    @{ Paragraph = \"Category cannot be null.\"; Highlight = \"cannot be null.\" }\n)\n\nforeach ($target in $targets) {\n    $wanted = $target.Paragraph + [char]13\n    $match = $null\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($p.Range.Text -eq $wanted) {\n            $match = $p\n            break\n        }\n    }\n    if ($match -eq $null) {\n        continue\n    }\n\n    $pr = $match.Range\n    $pr.Find.ClearFormatting()\n    $pr.Find.MatchCase = $true\n    $found = $pr.Find.Execute($target.Highlight)\n    if ($found) {\n        # Assigning on .Font (rather than the range shorthand) highlights\n        # only the matched sub-range, causing Word to split the original\n        # run the same way the UI would.\n        $pr.Font.HighlightColorIndex = [Microsoft.Office.Interop.Word.WdColorIndex]::wdYellow\n    }\n}\n"}
